$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet from "!_Schema" to "!!_Schema"
$ws.Name = "!!_Schema"

# Update the schema header/version marker cell
$ws.Range("A1").Value = "!!ObjTables type='Schema' tableFormat='row'"

# Rename "Model" -> "Class" for the Transaction row's Type column
$ws.Range("B3").Value = "Class"

# Update the view: zoom 160 -> 120, and move the selection from D5 to B4
$win = $excel.ActiveWindow
$win.Zoom = 120
$ws.Range("B4").Select()
